# Insert a new weekly record row above row 483, shifting all subsequent
# rows down by one (the last existing row becomes row 604 unchanged),
# then populate the newly inserted row 483 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 483; existing rows 483:603 shift to 484:604
$ws.Rows(483).Insert()

# Populate the new row 483 with the new data record
$ws.Range("A483").Value = 9
$ws.Range("B483").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C483").Value = "Metropolitana"
$ws.Range("D483").Value = 45135
$ws.Range("E483").Value = 13
$ws.Range("F483").Value = 100112039
$ws.Range("G483").Value = "Ciboulette"
$ws.Range("H483").Value = "Sin especificar"
$ws.Range("I483").Value = "Primera"
$ws.Range("J483").Value = 340
$ws.Range("K483").Value = 1000
$ws.Range("L483").Value = 1200
$ws.Range("M483").Value = 1100
$ws.Range("N483").Value = "`$/docena de atados"
$ws.Range("O483").Value = "Región Metropolitana"
$ws.Range("P483").Value = 367
$ws.Range("Q483").Value = 3
$ws.Range("R483").Value = "Hortaliza"
